$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (Wins / Losses / Ties) in row 1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting of the existing header cell (AC1) onto the new
# header cells so they pick up the same bold/centered/bordered style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the season record (Wins/Losses/Ties) for every data row.
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 30).Value = 93
    $ws.Cells.Item($r, 31).Value = 69
    $ws.Cells.Item($r, 32).Value = 0
}
